# Updated cryptos list on Sun Sep 17 16:26:33 UTC 2023 with GitHub Actions
#
# Applies the per-row Price (D) / Volume(1h) (E) updates for rows 2-48,
# and reflects a new "BabyDogeCoin" entry that is newly inserted above the
# former EnergySwap/Algorand/Mantle tail, shifting EnergySwap -> row 50,
# Algorand -> row 51 and dropping Mantle off the bottom of the (fixed
# size, 50-row) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-48: update Price (D) and Volume(1h) (E) values in place ---
$ws.Range("D2").Value = '26.740.18'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.639.97'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.63'
$ws.Range("E5").Value = '  +1.14%  '

$ws.Range("E7").Value = '  +0.39%  '

$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0624'
$ws.Range("E9").Value = '  -0.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.06'
$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").Value = '1.868.06'
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("D13").Value = '1.653.95'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("E14").Value = '  -0.58%  '

$ws.Range("E15").Value = '  -0.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.60'
$ws.Range("E16").Value = '  -0.40%  '

$ws.Range("D17").Value = '26.732.32'
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("E18").Value = '  -1.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '212.72'
$ws.Range("E19").Value = '  -1.48%  '

$ws.Range("E20").Value = '  +0.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.35'
$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.20'
$ws.Range("E22").Value = '  -0.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.32'
$ws.Range("E23").Value = '  +4.16%  '

$ws.Range("E24").Value = '  -2.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.46'
$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("E26").Value = '  +0.24%  '

$ws.Range("E27").Value = '  -1.48%  '

$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.61'
$ws.Range("E29").Value = '  -0.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0507'
$ws.Range("E30").Value = '  -1.09%  '

$ws.Range("E31").Value = '  +1.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.38'
$ws.Range("E32").Value = '  +1.10%  '

$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("D34").Value = '1.279.17'
$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("E35").Value = '  -0.55%  '

$ws.Range("E36").Value = '  +0.82%  '

$ws.Range("E37").Value = '  -1.17%  '

$ws.Range("E38").Value = '  +0.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.813'
$ws.Range("E39").Value = '  -0.77%  '

$ws.Range("E40").Value = '  +0.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.805'
$ws.Range("E41").Value = '  -0.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.21'
$ws.Range("E42").Value = '  -1.66%  '

$ws.Range("D43").Value = '1.777.98'
$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.25'
$ws.Range("E44").Value = '  -3.07%  '

$ws.Range("E45").Value = '  +3.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.24'
$ws.Range("E46").Value = '  -0.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.58'
$ws.Range("E47").Value = '  -1.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0522'
$ws.Range("E48").Value = '  +1.50%  '

# --- Rows 49-51: BabyDogeCoin is newly listed at row 49, pushing
#     EnergySwap down to row 50 and Algorand down to row 51; Mantle
#     (formerly row 51) falls off the bottom of the fixed-length table. ---

# Row 49 was EnergySwap -> now BabyDogeCoin
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$subsix = [string][char]0x2086
$ws.Range("D49").Value = '0.0' + $subsix + '0102'
$ws.Range("E49").Value = '  -2.81%  '

# Row 50 was Algorand -> now EnergySwap
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.60'
$ws.Range("E50").Value = '  -1.70%  '

# Row 51 was Mantle -> now Algorand
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0963'
$ws.Range("E51").Value = '  -0.05%  '
